$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1_MATERIALS_TURBINE")

# Switch to the turbine sheet and select the cell the user had clicked
# before removing the row (matches the post-edit active-cell state).
$ws.Activate()
$ws.Range("C10").Select()

# The "2000kW" row (row 6) only ever held a blank placeholder for an
# unused offshore 2MW turbine entry - remove it entirely, shifting the
# remaining turbine rows up.
$ws.Rows(6).Delete() | Out-Null
